$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 90
$linkCols = @("S","T","V","W","X","Y")

for ($row = 2; $row -le $lastRow; $row++) {
    # Column C: "Förändrad" date changed from 45184 to 45186
    $cDate = $ws.Range("C" + $row)
    if ($cDate.Value2() -eq 45184) {
        $cDate.Value2 = 45186
    }

    # Columns S,T,V,W,X,Y: append the "Beteckning" (column A) as the
    # second HYPERLINK() argument, turning it into the link's friendly text.
    $name = $ws.Range("A" + $row).Value()
    foreach ($col in $linkCols) {
        $cell = $ws.Range($col + $row)
        $f = $cell.Formula()
        if ($f -and $f -match '^=HYPERLINK\("[^"]*"\)$') {
            $newFormula = $f -replace '\)$', (', "' + $name + '")')
            $cell.Formula = $newFormula
        }
    }
}
